# Auto-generated Excel COM-interop edit script
# Implements: sheet renames + NAV/TER/rating-allocation data refresh

$wb = $excel.ActiveWorkbook

# ---- 1. Rename worksheets to reflect new fund naming convention ----
$wb.Worksheets.Item(2).Name = "Dynamic Bond"
$wb.Worksheets.Item(3).Name = "Corporate Bond Fund"
$wb.Worksheets.Item(4).Name = "Other Bond"
$wb.Worksheets.Item(7).Name = "Banking and PSU Fund"
$wb.Worksheets.Item(13).Name = "Ultra Short Duration Fund"
$wb.Worksheets.Item(15).Name = "Medium to Long Duration Fund"
$wb.Worksheets.Item(16).Name = "Gilt Fund with 10 year constant duration"

# ---- 2. Update cell values (TER / volatility / returns / rating allocations) ----
# Sheet 1: Liquid Fund
$ws = $wb.Worksheets.Item(1)
$cellValues = @(
    "G4=N:4.59",
    "H4=N:6.37",
    "I4=N:6.74",
    "H5=N:6.25",
    "F9=N:0.06",
    "H9=N:6.12",
    "I9=N:6.56",
    "F12=N:0.04",
    "H12=N:6.13",
    "I12=N:6.56",
    "H15=N:6.04",
    "I15=N:6.51",
    "H16=N:6.16",
    "H17=N:6.3",
    "G18=N:4.3",
    "H18=N:6.13",
    "I18=N:6.59",
    "I20=N:6.69",
    "F21=N:0.07",
    "G21=N:4.39",
    "H21=N:6.26",
    "G22=N:4.58",
    "H22=N:6.27",
    "F23=N:0.06",
    "H23=N:6.19",
    "I23=N:6.6",
    "K23=N:38.54",
    "L23=N:35.1",
    "M23=N:5.21",
    "F25=N:0.05",
    "G25=N:3.62",
    "H25=N:5.52",
    "I25=N:6.02",
    "H28=N:6.31",
    "I28=N:6.52",
    "L29=N:10.91",
    "M29=N:68.3",
    "G31=N:4.57",
    "H31=N:6.25",
    "I31=N:6.67",
    "F32=N:0.07",
    "H32=N:6.2",
    "I32=N:6.54",
    "F33=N:0.05",
    "H33=N:6.23",
    "I33=N:6.63",
    "H34=N:6.26",
    "I34=N:6.67",
    "H35=N:6.21",
    "H36=N:6.18",
    "F37=N:0.05",
    "G37=N:4.39",
    "H37=N:6.25",
    "I37=N:6.7",
    "F38=N:0.06",
    "G38=N:3.84",
    "H38=N:5.93",
    "I38=N:6.53",
    "F40=N:0.06",
    "G40=N:3.79",
    "H40=N:5.83",
    "F42=N:0.06",
    "G43=N:3.44",
    "F44=N:0.1",
    "G44=N:5.42"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 2: Dynamic Bond
$ws = $wb.Worksheets.Item(2)
$cellValues = @(
    "K16=N:77.77",
    "M16=N:7.86",
    "K23=N:90.27",
    "M23=N:4.42"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 3: Corporate Bond Fund
$ws = $wb.Worksheets.Item(3)
$cellValues = @(
    "K15=N:29.09",
    "M15=N:63.43"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 6: Gilt Fund
$ws = $wb.Worksheets.Item(6)
$cellValues = @(
    "K20=N:93.32",
    "K21=N:93.32",
    "K22=N:93.32"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 10: Short Duration Fund
$ws = $wb.Worksheets.Item(10)
$cellValues = @(
    "K16=N:23.06",
    "M16=N:66.86",
    "Q16=N:4.11",
    "K19=N:34.5",
    "M19=N:59.24"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 12: Low Duration Fund
$ws = $wb.Worksheets.Item(12)
$cellValues = @(
    "K18=N:13.47",
    "L18=N:16.75",
    "M18=N:60.18",
    "K21=N:30.32",
    "L21=N:4.97",
    "M21=N:37.76"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 13: Ultra Short Duration Fund
$ws = $wb.Worksheets.Item(13)
$cellValues = @(
    "K18=N:57.63",
    "L18=N:3.76",
    "M18=N:23.72",
    "K31=N:30.87",
    "M31=N:47.15"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 15: Medium to Long Duration Fund
$ws = $wb.Worksheets.Item(15)
$cellValues = @(
    "K13=N:85.79",
    "M13=S:-",
    "K14=N:62.54",
    "M14=N:32.93"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# Sheet 17: Overnight Fund
$ws = $wb.Worksheets.Item(17)
$cellValues = @(
    "F2=N:0.07",
    "G2=N:3.46",
    "I2=N:5.57",
    "G3=N:3.49",
    "H3=N:5.16",
    "I3=N:5.67",
    "F4=N:0.06",
    "G4=N:3.55",
    "H4=N:4.95",
    "G5=N:3.53",
    "H5=N:5.22",
    "F6=N:0.07",
    "G6=N:3.48",
    "F7=N:0.07",
    "G7=N:3.57",
    "G8=N:3.51",
    "F9=N:0.07",
    "G9=N:3.67",
    "F10=N:0.07",
    "G10=N:3.49",
    "G11=N:3.49",
    "F12=N:0.06",
    "G12=N:3.54",
    "G13=N:3.55",
    "F14=N:0.06",
    "G14=N:3.52",
    "F16=N:0.07",
    "G16=N:3.56",
    "G17=N:3.58",
    "F18=N:0.07",
    "G18=N:3.38",
    "F19=N:0.07",
    "G19=N:3.57",
    "G21=N:3.53",
    "G24=N:3.69",
    "F25=N:0.07",
    "G25=N:3.38",
    "G26=N:3.46",
    "F27=N:0.07",
    "G27=N:3.67",
    "G28=N:3.73"
)
foreach ($entry in $cellValues) {
    $parts = $entry.Split("=")
    $ref = $parts[0]
    $rest = $parts[1]
    $tag = $rest.Substring(0,2)
    $val = $rest.Substring(2)
    if ($tag -eq "N:") {
        $ws.Range($ref).Value = [double]$val
    } else {
        $ws.Range($ref).Value = $val
    }
}

# ---- 3. Re-color rating-allocation cells whose magnitude band changed ----
# Sheet 1: Liquid Fund
$ws = $wb.Worksheets.Item(1)
$ws.Range("K23").Font.Color = 15855596
$ws.Range("K23").Interior.Color = 6056192
$ws.Range("L23").Font.Color = 15855596
$ws.Range("L23").Interior.Color = 7043328
$ws.Range("L29").Font.Color = 0
$ws.Range("L29").Interior.Color = 12897152

# Sheet 3: Corporate Bond Fund
$ws = $wb.Worksheets.Item(3)
$ws.Range("K15").Font.Color = 0
$ws.Range("K15").Interior.Color = 8096000

# Sheet 10: Short Duration Fund
$ws = $wb.Worksheets.Item(10)
$ws.Range("K16").Font.Color = 0
$ws.Range("K16").Interior.Color = 8951296
$ws.Range("K19").Font.Color = 15855596
$ws.Range("K19").Interior.Color = 7043328

# Sheet 12: Low Duration Fund
$ws = $wb.Worksheets.Item(12)
$ws.Range("K18").Font.Color = 0
$ws.Range("K18").Interior.Color = 11318861
$ws.Range("L18").Font.Color = 0
$ws.Range("L18").Interior.Color = 11318861
$ws.Range("K21").Font.Color = 0
$ws.Range("K21").Interior.Color = 8096000
$ws.Range("L21").Font.Color = 0
$ws.Range("L21").Interior.Color = 14409650
$ws.Range("M21").Font.Color = 15855596
$ws.Range("M21").Interior.Color = 7043328

# Sheet 13: Ultra Short Duration Fund
$ws = $wb.Worksheets.Item(13)
$ws.Range("K18").Font.Color = 15855596
$ws.Range("K18").Interior.Color = 4214016
$ws.Range("L18").Font.Color = 0
$ws.Range("L18").Interior.Color = 14409650
$ws.Range("M18").Font.Color = 0
$ws.Range("M18").Interior.Color = 8951296
$ws.Range("K31").Font.Color = 0
$ws.Range("K31").Interior.Color = 8096000
$ws.Range("M31").Font.Color = 0
$ws.Range("M31").Interior.Color = 11318861

# Sheet 15: Medium to Long Duration Fund
$ws = $wb.Worksheets.Item(15)
$ws.Range("M13").Font.Color = 0
$ws.Range("M13").Interior.Color = 12959408
$ws.Range("M14").Font.Color = 0
$ws.Range("M14").Interior.Color = 8096000

